$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.88"
$ws.Range("E2").Value = "'0.91%"
$ws.Range("D3").Value = "'29.44"
$ws.Range("E3").Value = "'8.64%"
$ws.Range("D4").Value = "'5.166"
$ws.Range("E4").Value = "'1.94%"
$ws.Range("D5").Value = "'0.05713"
$ws.Range("E5").Value = "'0.45%"
$ws.Range("E6").Value = "'1.96%"
$ws.Range("D7").Value = "'0.8580"
$ws.Range("E7").Value = "'4.59%"
$ws.Range("D8").Value = "'0.8736"
$ws.Range("E8").Value = "'4.27%"
$ws.Range("D9").Value = "'0.1369"
$ws.Range("E9").Value = "'3.00%"
$ws.Range("D10").Value = "'0.07075"
$ws.Range("E10").Value = "'2.47%"
$ws.Range("D11").Value = "'0.02868"
$ws.Range("E11").Value = "'0.33%"
$ws.Range("D12").Value = "'0.09381"
$ws.Range("E12").Value = "'-0.20%"
$ws.Range("D13").Value = "'0.001515"
$ws.Range("E13").Value = "'-0.25%"
$ws.Range("D14").Value = "'0.04169"
$ws.Range("E14").Value = "'1.34%"
$ws.Range("D15").Value = "'0.0005987"
$ws.Range("E15").Value = "'0.16%"
$ws.Range("D16").Value = "'0.006187"
$ws.Range("E16").Value = "'0.83%"
$ws.Range("E17").Value = "'3,767.52%"
$ws.Range("E18").Value = "'-0.74%"
$ws.Range("D19").Value = "'3.060"
$ws.Range("E19").Value = "'1.98%"
$ws.Range("D20").Value = "'2.279"
$ws.Range("E20").Value = "'2.37%"
$ws.Range("D21").Value = "'0.3172"
$ws.Range("E21").Value = "'0.73%"
$ws.Range("D22").Value = "'0.03307"
$ws.Range("E22").Value = "'3.67%"
$ws.Range("E23").Value = "'0.39%"
$ws.Range("D24").Value = "'3.463"
$ws.Range("E24").Value = "'-3.08%"
$ws.Range("D25").Value = "'0.1379"
$ws.Range("E25").Value = "'0.44%"
$ws.Range("E26").Value = "'27.61%"
$ws.Range("D27").Value = "'0.001219"
$ws.Range("E27").Value = "'0.07%"
$ws.Range("D28").Value = "'0.0001209"
$ws.Range("E28").Value = "'23.42%"
$ws.Range("D40").Value = "'0.03756"
$ws.Range("E40").Value = "'1.62%"
$ws.Range("D41").Value = "'0.005787"
$ws.Range("E41").Value = "'68.86%"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("E42").Value = "'1.67%"
$ws.Range("D43").Value = "'0.002099"
$ws.Range("E43").Value = "'-10.67%"
$ws.Range("D44").Value = "'0.01025"
$ws.Range("E44").Value = "'9.32%"
$ws.Range("D45").Value = "'0.00005168"
$ws.Range("E45").Value = "'-0.58%"
$ws.Range("E46").Value = "'0.00%"
$ws.Range("D47").Value = "'0.07096"
$ws.Range("E47").Value = "'-30.05%"
$ws.Range("D48").Value = "'0.002571"
$ws.Range("E48").Value = "'-0.71%"
$ws.Range("E49").Value = "'0.00%"
$ws.Range("E50").Value = "'0.00%"
